$wb = $excel.ActiveWorkbook

# OFF sheet (sheet1) - row 2 ("H") updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 266
$wsOff.Range("C2").Value = 172
$wsOff.Range("D2").Value = 71
$wsOff.Range("E2").Value = 28

# DEF sheet (sheet2) - row 2 ("H") updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 370
$wsDef.Range("C2").Value = 267
$wsDef.Range("D2").Value = 79
$wsDef.Range("E2").Value = 42
